# Fix: add new properties at promocion and compra
#
# Block 1 - "Promoción" sheet section (rows 7-10): a new field
# "cantUnidadesDisponibles" is inserted into the small
# idProducto/idContenedorActual/existenciasActuales/nivelReorden relation
# table. The table (previously columns I-L) shifts one column right to
# J-M, leaving column I empty, and the new field lands in column H.
#
# Block 2 - "Compra" sheet section (rows 17-20): a new field
# "cantUnidadesCompradas" is inserted into the Productos_Compra relation
# table (replacing the stray leftover "O" cell at J18), and the
# Pedidos_Proveedor mini-table (previously columns K-L) shifts one
# column right to L-M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE on ordering: new shared strings are appended to xl/sharedStrings.xml
# in the order the cell values are *written*, not in sheet (row/col) order.
# The target file has "cantUnidadesCompradas" (Compra block) ahead of
# "cantUnidadesDisponibles" (Promoción block) in the shared-string table,
# so the Compra block's new-string write must happen first.

# ---------------------------------------------------------------------
# Block 2: rows 17-20 ("Compra" / Productos_Compra / Pedidos_Proveedor)
# ---------------------------------------------------------------------

# Row 17: move the "Pedidos_Proveedor" label from K17 to L17.
$ws.Range("K17").Copy()
$ws.Range("L17").PasteSpecial(-4122)
$ws.Range("L17").Value = $ws.Range("K17").Value2
$ws.Range("K17").Clear()

# Row 18: shift L18->M18, K18->L18 (right to left), then replace J18
# (previously the stray "O" cell) with the new
# "cantUnidadesCompradas" header, matching I18's header style.
$ws.Range("L18").Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M18").Value = $ws.Range("L18").Value2

$ws.Range("K18").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L18").Value = $ws.Range("K18").Value2
$ws.Range("K18").Clear()

$ws.Range("I18").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("J18").Value = "cantUnidadesCompradas"

# Row 19: shift L19->M19, K19->L19 (right to left), then place the new
# "NN" type cell in J19, matching I19's style.
$ws.Range("L19").Copy()
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("M19").Value = $ws.Range("L19").Value2

$ws.Range("K19").Copy()
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("L19").Value = $ws.Range("K19").Value2
$ws.Range("K19").Clear()

$ws.Range("I19").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("J19").Value = "NN"

# Row 20: blank divider row - shift the formatting the same way.
$ws.Range("L20").Copy()
$ws.Range("M20").PasteSpecial(-4122)

$ws.Range("K20").Copy()
$ws.Range("L20").PasteSpecial(-4122)
$ws.Range("K20").Clear()

$ws.Range("I20").Copy()
$ws.Range("J20").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Block 1: rows 7-10 ("Promoción" / idProducto-Contenedor relation)
# ---------------------------------------------------------------------

# Row 7: move the "idProducto" label from I7 to J7.
$ws.Range("I7").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J7").Value = $ws.Range("I7").Value2
$ws.Range("I7").Clear()

# Row 8: shift L8->M8, K8->L8, J8->K8, I8->J8 (right to left), then
# place the new "cantUnidadesDisponibles" header in H8.
$ws.Range("L8").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("M8").Value = $ws.Range("L8").Value2

$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").Value = $ws.Range("K8").Value2

$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = $ws.Range("J8").Value2

$ws.Range("I8").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Value = $ws.Range("I8").Value2

$ws.Range("I8").Clear()

$ws.Range("G8").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Value = "cantUnidadesDisponibles"

# Row 9: shift L9->M9, K9->L9, J9->K9, I9->J9 (right to left), then
# place the new "NN" type cell in H9.
$ws.Range("L9").Copy()
$ws.Range("M9").PasteSpecial(-4122)
$ws.Range("M9").Value = $ws.Range("L9").Value2

$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = $ws.Range("K9").Value2

$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Value = $ws.Range("J9").Value2

$ws.Range("I9").Copy()
$ws.Range("J9").PasteSpecial(-4122)
$ws.Range("J9").Value = $ws.Range("I9").Value2

$ws.Range("I9").Clear()

$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = "NN"

# Row 10: blank divider row - shift the formatting the same way.
$ws.Range("L10").Copy()
$ws.Range("M10").PasteSpecial(-4122)

$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial(-4122)

$ws.Range("I10").Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("I10").Clear()

$ws.Range("G10").Copy()
$ws.Range("H10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet view / window bookkeeping to match the saved session state.
# ---------------------------------------------------------------------
$ws.Range("K16").Select()
$excel.ActiveWindow.ScrollRow = 1
